# Update column G ("K") values on Sheet1 per the recalculated save_data.
# Mapping of row -> new value for column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 4
    16 = 2
    17 = 2
    18 = 3
    19 = 1
    20 = 2
    21 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
